# Generate Report for Handoff
#
# File "8de05a06-e841-430d-ad60-ba25ede17482.md" has just been handed off
# again (a new Xliff/handoff generation cycle). Update the three
# timestamp cells that track this event across the workbook's sheets:
#
#   - Overview!G4  "Latest HO Xliff Generate Date" -> 2016-10-19 16:10:00
#   - zh-cn!H4      "Latest Handoff Datetime"       -> 2016-10-19 16:09:49
#   - de-de!H4      "Latest Handoff Datetime"       -> 2016-10-19 16:10:00

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-10-19 16:10:00"
$wsZhCn.Range("H4").Value     = "2016-10-19 16:09:49"
$wsDeDe.Range("H4").Value     = "2016-10-19 16:10:00"
